# "Updating new data files" — Sheet1's A2/A3 labels were shortened
# (S1 -> S, B1 -> B) and the live selection moved to A4.
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("A2").Value = "S"
$ws1.Range("A3").Value = "B"

$ws1.Range("A4").Select() | Out-Null
